$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1352.6216
$ws.Range("I15").Value = 1352.6216
$ws.Range("K15").Value = 4057.8648
$ws.Range("M15").Value = -3888.8648
$ws.Range("H40").Value = 166669140
$ws.Range("I40").Value = 3500
$ws.Range("K40").Value = 3500
$ws.Range("M40").Value = -3325
$ws.Range("H92").Value = 1915.1578
$ws.Range("I92").Value = 795.8
$ws.Range("K92").Value = 795.8
$ws.Range("M92").Value = 452.2
$ws.Range("H99").Value = 2273.125
$ws.Range("I99").Value = 162.66667
$ws.Range("J99").Value = 3539.4
$ws.Range("K99").Value = 488.00001
$ws.Range("L99").Value = 10618.2
$ws.Range("M99").Value = 1009.99999
$ws.Range("N99").Value = -13614.2
$ws.Range("H101").Value = 399
$ws.Range("J101").Value = 449.9091
$ws.Range("L101").Value = 1349.7273
$ws.Range("N101").Value = -4593.7273
$ws.Range("H107").Value = 846.2222
$ws.Range("J107").Value = 1821
$ws.Range("L107").Value = 1821
$ws.Range("N107").Value = -5661
$ws.Range("H127").Value = 2331
$ws.Range("I127").Value = 1001.5
$ws.Range("K127").Value = 3004.5
$ws.Range("M127").Value = 1955.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 50166.668
$ws.Range("J30").Value = 75000
$ws.Range("L30").Value = 75000
$ws.Range("N30").Value = -75300
$ws.Range("H32").Value = 3257.6516
$ws.Range("I32").Value = 2384.6206
$ws.Range("K32").Value = 2384.6206
$ws.Range("M32").Value = -2097.6206
$ws.Range("H74").Value = 2644.0952
$ws.Range("I74").Value = 1515.1333
$ws.Range("K74").Value = 1515.1333
$ws.Range("M74").Value = -641.1333
$ws.Range("H77").Value = 2644.0952
$ws.Range("I77").Value = 1515.1333
$ws.Range("K77").Value = 7575.666499999999
$ws.Range("M77").Value = -3207.666499999999
$ws.Range("H106").Value = 55000
$ws.Range("J106").Value = 55000
$ws.Range("L106").Value = 55000
$ws.Range("N106").Value = -57524
$ws.Range("H122").Value = 5848.857
$ws.Range("I122").Value = 5925.769
$ws.Range("K122").Value = 17777.307
$ws.Range("M122").Value = -15327.307
$ws.Range("H132").Value = 2783416
$ws.Range("I132").Value = 5334.483
$ws.Range("K132").Value = 16003.449
$ws.Range("M132").Value = -13473.449

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2289.8635
$ws.Range("I94").Value = 2317.9
$ws.Range("K94").Value = 2317.9
$ws.Range("M94").Value = -1866.9

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31253134
$ws.Range("I31").Value = 47621916
$ws.Range("J31").Value = 3643.0908
$ws.Range("K31").Value = 47621916
$ws.Range("L31").Value = 3643.0908
$ws.Range("M31").Value = -47621621
$ws.Range("N31").Value = -4233.0908
$ws.Range("H34").Value = 31253134
$ws.Range("I34").Value = 47621916
$ws.Range("J34").Value = 3643.0908
$ws.Range("K34").Value = 47621916
$ws.Range("L34").Value = 3643.0908
$ws.Range("M34").Value = -47621714
$ws.Range("N34").Value = -4047.0908
$ws.Range("H94").Value = 1000.82355
$ws.Range("I94").Value = 933.875
$ws.Range("K94").Value = 933.875
$ws.Range("M94").Value = -482.875
$ws.Range("H122").Value = 4216.5557
$ws.Range("I122").Value = 3601.8
$ws.Range("K122").Value = 10805.4
$ws.Range("M122").Value = -8355.400000000001
$ws.Range("H132").Value = 2877.5833
$ws.Range("I132").Value = 2807.389
$ws.Range("J132").Value = 3088.1667
$ws.Range("K132").Value = 8422.167000000001
$ws.Range("L132").Value = 9264.500100000001
$ws.Range("M132").Value = -5892.167000000001
$ws.Range("N132").Value = -14324.5001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 26592.555
$ws.Range("I76").Value = 9000
$ws.Range("J76").Value = 28791.625
$ws.Range("K76").Value = 27000
$ws.Range("L76").Value = 86374.875
$ws.Range("M76").Value = -26617
$ws.Range("N76").Value = -87140.875
$ws.Range("H79").Value = 26592.555
$ws.Range("I79").Value = 9000
$ws.Range("J79").Value = 28791.625
$ws.Range("K79").Value = 27000
$ws.Range("L79").Value = 86374.875
$ws.Range("M79").Value = -25674
$ws.Range("N79").Value = -89026.875
$ws.Range("H86").Value = 637.875
$ws.Range("J86").Value = 980.6
$ws.Range("L86").Value = 2941.8
$ws.Range("N86").Value = -5313.8
$ws.Range("H87").Value = 9157
$ws.Range("I87").Value = 1098.3334
$ws.Range("K87").Value = 3295.0002
$ws.Range("M87").Value = -2047.0002
$ws.Range("H89").Value = 637.875
$ws.Range("J89").Value = 980.6
$ws.Range("L89").Value = 8825.4
$ws.Range("N89").Value = -20681.4
$ws.Range("H90").Value = 9157
$ws.Range("I90").Value = 1098.3334
$ws.Range("K90").Value = 9885.000599999999
$ws.Range("M90").Value = -3645.000599999999
$ws.Range("H98").Value = 750
$ws.Range("I98").Value = 750
$ws.Range("K98").Value = 2250
$ws.Range("M98").Value = -752
$ws.Range("H129").Value = 6265.769
$ws.Range("I129").Value = 3262.5
$ws.Range("J129").Value = 11071
$ws.Range("K129").Value = 9787.5
$ws.Range("L129").Value = 33213
$ws.Range("M129").Value = -4787.5
$ws.Range("N129").Value = -43213
$ws.Range("H137").Value = 6966.5884
$ws.Range("I137").Value = 1731.1
$ws.Range("J137").Value = 14445.857
$ws.Range("K137").Value = 5193.299999999999
$ws.Range("L137").Value = 43337.571
$ws.Range("M137").Value = -93.29999999999927
$ws.Range("N137").Value = -53537.571

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 130899
$ws.Range("J39").Value = 130899
$ws.Range("L39").Value = 130899
$ws.Range("N39").Value = -131963
$ws.Range("H97").Value = 1121.8889
$ws.Range("I97").Value = 1299.7142
$ws.Range("J97").Value = 499.5
$ws.Range("K97").Value = 1299.7142
$ws.Range("L97").Value = 499.5
$ws.Range("M97").Value = -803.7141999999999
$ws.Range("N97").Value = -1491.5
$ws.Range("H134").Value = 89569.25
$ws.Range("J134").Value = 89569.25
$ws.Range("L134").Value = 268707.75
$ws.Range("N134").Value = -273777.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1831
$ws.Range("I46").Value = 1466.3334
$ws.Range("K46").Value = 1466.3334
$ws.Range("M46").Value = -1278.3334
$ws.Range("H68").Value = 8335633
$ws.Range("I68").Value = 10418541
$ws.Range("J68").Value = 4003
$ws.Range("K68").Value = 10418541
$ws.Range("L68").Value = 4003
$ws.Range("M68").Value = -10417792
$ws.Range("N68").Value = -5501
$ws.Range("H71").Value = 8335633
$ws.Range("I71").Value = 10418541
$ws.Range("J71").Value = 4003
$ws.Range("K71").Value = 52092705
$ws.Range("L71").Value = 20015
$ws.Range("M71").Value = -52088961
$ws.Range("N71").Value = -27503

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 99998.5
$ws.Range("J95").Value = 99998.5
$ws.Range("L95").Value = 99998.5
$ws.Range("N95").Value = -105490.5
$ws.Range("H100").Value = 777.06665
$ws.Range("I100").Value = 777.06665
$ws.Range("K100").Value = 1554.1333
$ws.Range("M100").Value = -1013.1333
$ws.Range("H132").Value = 425020.97
$ws.Range("I132").Value = 7974.1055
$ws.Range("K132").Value = 23922.3165
$ws.Range("M132").Value = -21392.3165
